$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Post content" bullet is removed, and the "_GoBack" bookmark (which
# previously sat, collapsed/empty, on the bullet right after "Name
# Ideas") is relocated onto the start of the "Extended description
# (Optional)" bullet, whose text is then trimmed to just
# "description (Optional)".
# ------------------------------------------------------------------

# Step 1: find and delete the "Post content" paragraph.
$postContent = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r`a") -eq "Post content") {
        $postContent = $candidate
        break
    }
}
if ($postContent -ne $null) {
    [void]$postContent.Range.Delete()
}

# Step 2: locate the "Extended description (Optional)" paragraph and
# re-plant the _GoBack bookmark at its very start. Bookmarks.Add with
# a name that already exists elsewhere moves the bookmark rather than
# duplicating it, so this both removes it from its old (empty) home
# and drops a collapsed bookmark right before the run here.
$descPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("Extended description (Optional)")) {
        $descPara = $candidate
        break
    }
}

if ($descPara -ne $null) {
    $startPos = $descPara.Range.Start
    $collapsed = $d.Range($startPos, $startPos)
    [void]$d.Bookmarks.Add("_GoBack", $collapsed)

    # Step 3: drop the leading "Extended " from the bullet text.
    [void]$descPara.Range.Find.Execute("Extended description (Optional)", $true, $false, $false, $false, $false,
                                        $true, 1, $false, "description (Optional)", 2)
}
